# #6665 - Adding unit for autodiag text responses
#
# Inserts a new "unite_reponse" header column between "libelle_question"
# (column F) and "format_reponse" (previously column G, now H) on the
# "questions" sheet, shifting the remaining headers one column to the
# right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("questions")

# New column goes in right before "format_reponse", which currently lives
# in column 7 (G). Capture the existing values from G1:L1 so they can be
# re-written one column further right (H1:M1).
$firstCol = 7
$lastCol = 12

$values = @()
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $values += $ws.Cells.Item(1, $c).Value()
}

# Write them back starting one column to the right, working from the end
# so no value is overwritten before it's been read.
for ($i = $values.Length - 1; $i -ge 0; $i--) {
    $ws.Cells.Item(1, $firstCol + $i + 1).Value = $values[$i]
}

# Populate the freed-up column with the new header label.
$ws.Cells.Item(1, $firstCol).Value = "unite_reponse"
